# Scheduled runner update: refresh market-price / profit figures across
# several Leve sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H41").Value  = 362.8
$ws.Range("I41").Value  = 344.33334
$ws.Range("J41").Value  = 365.31818
$ws.Range("K41").Value  = 344.33334
$ws.Range("L41").Value  = 365.31818
$ws.Range("M41").Value  = 95.66665999999998
$ws.Range("N41").Value  = -1245.31818

$ws.Range("H113").Value = 54396.105
$ws.Range("J113").Value = 1917.5834
$ws.Range("L113").Value = 1917.5834
$ws.Range("N113").Value = -8425.5834

$ws.Range("H129").Value = 1004.7931
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1004.7931
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3014.3793
$ws.Range("M129").Value = ""
$ws.Range("N129").Value = -13014.3793

$ws.Range("H137").Value = 1334.3334
$ws.Range("I137").Value = 1096.1945
$ws.Range("J137").Value = 2763.1667
$ws.Range("K137").Value = 3288.5835
$ws.Range("L137").Value = 8289.500100000001
$ws.Range("M137").Value = -738.5835000000002
$ws.Range("N137").Value = -13389.5001

$ws.Range("H138").Value = 5733.6865
$ws.Range("I138").Value = 1965.1111
$ws.Range("J138").Value = 6777.2925
$ws.Range("K138").Value = 5895.3333
$ws.Range("L138").Value = 20331.8775
$ws.Range("M138").Value = -755.3333000000002
$ws.Range("N138").Value = -30611.8775

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H19").Value  = 0
$ws.Range("I19").Value  = 0
$ws.Range("K19").Value  = 0
$ws.Range("M19").Value  = ""

$ws.Range("H110").Value = 62626144
$ws.Range("I110").Value = 100200960
$ws.Range("J110").Value = 1446
$ws.Range("K110").Value = 100200960
$ws.Range("L110").Value = 1446
$ws.Range("M110").Value = -100198915
$ws.Range("N110").Value = -5536

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H4").Value   = 3666.6667
$ws.Range("I4").Value   = 1500
$ws.Range("J4").Value   = 8000
$ws.Range("K4").Value   = 1500
$ws.Range("L4").Value   = 8000
$ws.Range("M4").Value   = -1388
$ws.Range("N4").Value   = -8224

$ws.Range("H31").Value  = 27167.145
$ws.Range("I31").Value  = 974.26086
$ws.Range("J31").Value  = 42614.23
$ws.Range("K31").Value  = 974.26086
$ws.Range("L31").Value  = 42614.23
$ws.Range("M31").Value  = -679.26086
$ws.Range("N31").Value  = -43204.23

$ws.Range("H34").Value  = 27167.145
$ws.Range("I34").Value  = 974.26086
$ws.Range("J34").Value  = 42614.23
$ws.Range("K34").Value  = 974.26086
$ws.Range("L34").Value  = 42614.23
$ws.Range("M34").Value  = -772.26086
$ws.Range("N34").Value  = -43018.23

$ws.Range("H110").Value = 38000
$ws.Range("J110").Value = 38000
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180

$ws.Range("H134").Value = 1217
$ws.Range("I134").Value = 1096.5
$ws.Range("J134").Value = 1638.75
$ws.Range("K134").Value = 3289.5
$ws.Range("L134").Value = 4916.25
$ws.Range("M134").Value = -754.5
$ws.Range("N134").Value = -9986.25

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H3").Value   = 6312
$ws.Range("I3").Value   = 6312
$ws.Range("K3").Value   = 18936
$ws.Range("M3").Value   = -18824

$ws.Range("H4").Value   = 111111250
$ws.Range("I4").Value   = 150
$ws.Range("K4").Value   = 450
$ws.Range("M4").Value   = -338

$ws.Range("H124").Value = 3465
$ws.Range("I124").Value = 2015
$ws.Range("J124").Value = 4190
$ws.Range("K124").Value = 6045
$ws.Range("L124").Value = 12570
$ws.Range("M124").Value = -1135
$ws.Range("N124").Value = -22390

$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = ""

$ws.Range("H131").Value = 844.53
$ws.Range("J131").Value = 860.0319
$ws.Range("L131").Value = 2580.0957
$ws.Range("N131").Value = -12660.0957

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H34").Value  = 0
$ws.Range("I34").Value  = 0
$ws.Range("K34").Value  = 0
$ws.Range("M34").Value  = ""

$ws.Range("H76").Value  = 0
$ws.Range("I76").Value  = 0
$ws.Range("K76").Value  = 0
$ws.Range("M76").Value  = ""

$ws.Range("H79").Value  = 0
$ws.Range("I79").Value  = 0
$ws.Range("K79").Value  = 0
$ws.Range("M79").Value  = ""

$ws.Range("H132").Value = 1842.4231
$ws.Range("I132").Value = 1430.15
$ws.Range("J132").Value = 3216.6667
$ws.Range("K132").Value = 4290.450000000001
$ws.Range("L132").Value = 9650.000100000001
$ws.Range("M132").Value = -1760.450000000001
$ws.Range("N132").Value = -14710.0001

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value  = 854.15
$ws.Range("I22").Value  = 1149.75
$ws.Range("J22").Value  = 780.25
$ws.Range("K22").Value  = 1149.75
$ws.Range("L22").Value  = 780.25
$ws.Range("M22").Value  = -854.75
$ws.Range("N22").Value  = -1370.25

$ws.Range("H24").Value  = 14235.333
$ws.Range("J24").Value  = 16900
$ws.Range("L24").Value  = 16900
$ws.Range("N24").Value  = -17586

$ws.Range("H27").Value  = 854.15
$ws.Range("I27").Value  = 1149.75
$ws.Range("J27").Value  = 780.25
$ws.Range("K27").Value  = 1149.75
$ws.Range("L27").Value  = 780.25
$ws.Range("M27").Value  = -1042.75
$ws.Range("N27").Value  = -994.25

$ws.Range("H40").Value  = 2250
$ws.Range("I40").Value  = 2000
$ws.Range("J40").Value  = 2500
$ws.Range("K40").Value  = 2000
$ws.Range("L40").Value  = 2500
$ws.Range("M40").Value  = -1864
$ws.Range("N40").Value  = -2772

$ws.Range("H122").Value = 2100
$ws.Range("I122").Value = 2150
$ws.Range("K122").Value = 6450
$ws.Range("M122").Value = -4000

$ws.Range("H132").Value = 2324.568
$ws.Range("I132").Value = 2416.262
$ws.Range("J132").Value = 399
$ws.Range("K132").Value = 7248.786
$ws.Range("L132").Value = 1197
$ws.Range("M132").Value = -4718.786
$ws.Range("N132").Value = -6257

# ---------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H20").Value  = 8654.546
$ws.Range("I20").Value  = 0
$ws.Range("J20").Value  = 8654.546
$ws.Range("K20").Value  = 0
$ws.Range("L20").Value  = 8654.546
$ws.Range("M20").Value  = ""
$ws.Range("N20").Value  = -9134.546

Write-Host "Aegis_Profits sheets updated."
